# Tambah referensi manual: 118976
#
# 1) Rename the worksheet tab from "Data Referensi" to "Sheet1".
# 2) Re-style the header row: drop the white-on-blue look (bold white font
#    on a solid blue fill) in favor of a plain bold font with no fill,
#    and change the vertical alignment from centered to top.
# 3) The reference columns (A = "Kode BA", F = "No") were stored as
#    numbers; store them as text instead (same displayed values).
# 4) Append a new manual reference row (row 79) for Kementerian Keuangan /
#    satker 118976.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename the sheet -----------------------------------------------
$ws.Name = "Sheet1"

# --- 2) Re-style the header row (row 1) --------------------------------
$header = $ws.Range("A1:F1")
$header.Font.Bold = $true
$header.Font.Color = 0            # plain black text instead of white
$header.Interior.Pattern = -4142  # xlNone -> no fill instead of solid blue
$header.HorizontalAlignment = -4108  # xlCenter (unchanged)
$header.VerticalAlignment = -4160    # xlTop (was xlCenter)

# --- 3) Convert existing "Kode BA" (A) and "No" (F) columns to text ----
$lastRow = 78
for ($r = 2; $r -le $lastRow; $r++) {
    $cellA = $ws.Cells.Item($r, 1)
    $valA = [string]$cellA.Value2
    $cellA.NumberFormat = "@"
    $cellA.Value = $valA

    $cellF = $ws.Cells.Item($r, 6)
    $valF = [string]$cellF.Value2
    $cellF.NumberFormat = "@"
    $cellF.Value = $valF
}

# --- 4) Append the new manual reference row (row 79) -------------------
$newRow = 79

$cellA79 = $ws.Cells.Item($newRow, 1)
$cellA79.NumberFormat = "@"
$cellA79.Value = "015"

$ws.Cells.Item($newRow, 2).Value = "Kementerian Keuangan"

$cellC79 = $ws.Cells.Item($newRow, 3)
$cellC79.NumberFormat = "@"
$cellC79.Value = "118976"

$ws.Cells.Item($newRow, 4).Value = "TEST"
$ws.Cells.Item($newRow, 5).Value = "TESTT TESTTTT"
$ws.Cells.Item($newRow, 6).Value = 78
